# "fixed dates some more in xlsx"
#
# The "Follow up Data" sheet has a "date" column (column L) whose values
# were entered as text like "6.25.25" and need to become "6.25.2025"
# (and similarly for a handful of other dates in the same column).
#
# Setting Range.Value / .Value2 / .Formula directly with a date-looking
# string such as "6.25.2025" makes Excel "smart" parse it into a real
# date serial number (and re-stamp the cell with a date number format),
# which is NOT what happened in the source edit (there the cell stayed
# plain text in the shared-string table). To reproduce that faithfully we
# build the text via a formula that evaluates to a string ("=""6.25.2025""")
# and then do a Copy / PasteSpecial(xlPasteValues) of that cell onto
# itself, which commits the literal text value without Excel's
# type-inference kicking in and without touching the cell's style.

function Set-ExactText($Worksheet, $StartRow, $EndRow, $Col, $Text) {
    $escaped = $Text.Replace('"', '""')
    for ($r = $StartRow; $r -le $EndRow; $r++) {
        $cell = $Worksheet.Cells.Item($r, $Col)
        $cell.Formula = '="' + $escaped + '"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Follow up Data")
$dateCol = 12  # column L ("date")

Set-ExactText $ws 2   99  $dateCol "6.25.2025"
Set-ExactText $ws 101 190 $dateCol "7.2.2025"
Set-ExactText $ws 192 256 $dateCol "7.8.2025"
Set-ExactText $ws 258 303 $dateCol "7.17.2025"
Set-ExactText $ws 305 315 $dateCol "7.22.2025"
